# Update NATMI ligand-receptor pair metrics (Fn1-Itgb1) with refreshed TPM-derived
# expression values. Ligand/Receptor average & total expression values (columns
# G, H, M, N) come from newly computed per-cluster TPM averages; all of the
# specificity and edge-weight columns (I, J, O, P, Q, R, S, T) are downstream
# values recomputed from those updated expression numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = 29.20950566666667
$ws.Range("H2").Value = 87.628517
$ws.Range("I2").Value = 0.01829497698069002
$ws.Range("J2").Value = 0.01840828041918582
$ws.Range("M2").Value = 121.928739
$ws.Range("N2").Value = 365.786217
$ws.Range("O2").Value = 0.2282232151508951
$ws.Range("P2").Value = 0.2419720431319445
$ws.Range("Q2").Value = 3561.478192750021
$ws.Range("R2").Value = 32053.30373475018
$ws.Range("S2").Value = 0.004175338467644691
$ws.Range("T2").Value = 0.004454289223576161
$ws.Range("G3").Value = 29.20950566666667
$ws.Range("H3").Value = 87.628517
$ws.Range("I3").Value = 0.01829497698069002
$ws.Range("J3").Value = 0.01840828041918582
$ws.Range("O3").Value = 0.2768624053389947
$ws.Range("P3").Value = 0.2935413991166814
$ws.Range("Q3").Value = 4320.504460316203
$ws.Range("R3").Value = 38884.54014284583
$ws.Range("S3").Value = 0.005065191332495378
$ws.Range("T3").Value = 0.005403592389580016
$ws.Range("G4").Value = 29.20950566666667
$ws.Range("H4").Value = 87.628517
$ws.Range("I4").Value = 0.01829497698069002
$ws.Range("J4").Value = 0.01840828041918582
$ws.Range("M4").Value = 83.50496933333334
$ws.Range("N4").Value = 250.514908
$ws.Range("O4").Value = 0.1563025480180701
$ws.Range("P4").Value = 0.1657186665504434
$ws.Range("Q4").Value = 2439.138874936827
$ws.Range("R4").Value = 21952.24987443144
$ws.Range("S4").Value = 0.002859551518013788
$ws.Range("T4").Value = 0.003050595684554111
$ws.Range("G5").Value = 29.20950566666667
$ws.Range("H5").Value = 87.628517
$ws.Range("I5").Value = 0.01829497698069002
$ws.Range("J5").Value = 0.01840828041918582
$ws.Range("M5").Value = 91.06846250000001
$ws.Range("N5").Value = 182.136925
$ws.Range("O5").Value = 0.1704597085236707
$ws.Range("P5").Value = 0.1204857969594293
$ws.Range("Q5").Value = 2660.064771448371
$ws.Range("R5").Value = 15960.38862869023
$ws.Range("S5").Value = 0.003118556443575686
$ws.Range("T5").Value = 0.002217936336958261
$ws.Range("G6").Value = 29.20950566666667
$ws.Range("H6").Value = 87.628517
$ws.Range("I6").Value = 0.01829497698069002
$ws.Range("J6").Value = 0.01840828041918582
$ws.Range("M6").Value = 89.83562999999999
$ws.Range("N6").Value = 269.50689
$ws.Range("O6").Value = 0.1681521229683693
$ws.Range("P6").Value = 0.1782820942415013
$ws.Range("Q6").Value = 2624.05434355357
$ws.Range("R6").Value = 23616.48909198213
$ws.Range("S6").Value = 0.003076339218960474
$ws.Range("T6").Value = 0.003281866784517269
$ws.Range("I7").Value = 0.913374480506715
$ws.Range("J7").Value = 0.9190311407684336
$ws.Range("M7").Value = 121.928739
$ws.Range("N7").Value = 365.786217
$ws.Range("O7").Value = 0.2282232151508951
$ws.Range("P7").Value = 0.2419720431319445
$ws.Range("Q7").Value = 177806.3616900137
$ws.Range("R7").Value = 1600257.255210123
$ws.Range("S7").Value = 0.208453260578021
$ws.Range("T7").Value = 0.2223798428336196
$ws.Range("I8").Value = 0.913374480506715
$ws.Range("J8").Value = 0.9190311407684336
$ws.Range("O8").Value = 0.2768624053389947
$ws.Range("P8").Value = 0.2935413991166814
$ws.Range("R8").Value = 1941305.894519621
$ws.Range("S8").Value = 0.2528790556483438
$ws.Range("T8").Value = 0.2697736868929658
$ws.Range("I9").Value = 0.913374480506715
$ws.Range("J9").Value = 0.9190311407684336
$ws.Range("M9").Value = 83.50496933333334
$ws.Range("N9").Value = 250.514908
$ws.Range("O9").Value = 0.1563025480180701
$ws.Range("P9").Value = 0.1657186665504434
$ws.Range("Q9").Value = 121773.7089874781
$ws.Range("R9").Value = 1095963.380887303
$ws.Range("S9").Value = 0.1427627585978806
$ws.Range("T9").Value = 0.1523006151664776
$ws.Range("I10").Value = 0.913374480506715
$ws.Range("J10").Value = 0.9190311407684336
$ws.Range("M10").Value = 91.06846250000001
$ws.Range("N10").Value = 182.136925
$ws.Range("O10").Value = 0.1704597085236707
$ws.Range("P10").Value = 0.1204857969594293
$ws.Range("Q10").Value = 132803.4072576479
$ws.Range("R10").Value = 796820.4435458872
$ws.Range("S10").Value = 0.1556935477201338
$ws.Range("T10").Value = 0.1107301994260182
$ws.Range("I11").Value = 0.913374480506715
$ws.Range("J11").Value = 0.9190311407684336
$ws.Range("M11").Value = 89.83562999999999
$ws.Range("N11").Value = 269.50689
$ws.Range("O11").Value = 0.1681521229683693
$ws.Range("P11").Value = 0.1782820942415013
$ws.Range("Q11").Value = 131005.591064386
$ws.Range("R11").Value = 1179050.319579474
$ws.Range("S11").Value = 0.1535858579623356
$ws.Range("T11").Value = 0.1638467964493523
$ws.Range("G12").Value = 57.98602933333333
$ws.Range("H12").Value = 173.958088
$ws.Range("I12").Value = 0.03631876156896331
$ws.Range("J12").Value = 0.03654368891224535
$ws.Range("M12").Value = 121.928739
$ws.Range("N12").Value = 365.786217
$ws.Range("O12").Value = 0.2282232151508951
$ws.Range("P12").Value = 0.2419720431319445
$ws.Range("Q12").Value = 7070.163436230343
$ws.Range("R12").Value = 63631.47092607308
$ws.Range("S12").Value = 0.008288784535567573
$ws.Range("T12").Value = 0.008842551069674194
$ws.Range("G13").Value = 57.98602933333333
$ws.Range("H13").Value = 173.958088
$ws.Range("I13").Value = 0.03631876156896331
$ws.Range("J13").Value = 0.03654368891224535
$ws.Range("O13").Value = 0.2768624053389947
$ws.Range("P13").Value = 0.2935413991166814
$ws.Range("Q13").Value = 8576.964678200344
$ws.Range("R13").Value = 77192.6821038031
$ws.Range("S13").Value = 0.01005529968691662
$ws.Range("T13").Value = 0.01072708557218526
$ws.Range("G14").Value = 57.98602933333333
$ws.Range("H14").Value = 173.958088
$ws.Range("I14").Value = 0.03631876156896331
$ws.Range("J14").Value = 0.03654368891224535
$ws.Range("M14").Value = 83.50496933333334
$ws.Range("N14").Value = 250.514908
$ws.Range("O14").Value = 0.1563025480180701
$ws.Range("P14").Value = 0.1657186665504434
$ws.Range("Q14").Value = 4842.121601241767
$ws.Range("R14").Value = 43579.09441117589
$ws.Range("S14").Value = 0.005676714974089726
$ws.Range("T14").Value = 0.006055971397371521
$ws.Range("G15").Value = 57.98602933333333
$ws.Range("H15").Value = 173.958088
$ws.Range("I15").Value = 0.03631876156896331
$ws.Range("J15").Value = 0.03654368891224535
$ws.Range("M15").Value = 91.06846250000001
$ws.Range("N15").Value = 182.136925
$ws.Range("O15").Value = 0.1704597085236707
$ws.Range("P15").Value = 0.1204857969594293
$ws.Range("Q15").Value = 5280.698537866567
$ws.Range("R15").Value = 31684.1912271994
$ws.Range("S15").Value = 0.006190885510986181
$ws.Range("T15").Value = 0.004402995482429342
$ws.Range("G16").Value = 57.98602933333333
$ws.Range("H16").Value = 173.958088
$ws.Range("I16").Value = 0.03631876156896331
$ws.Range("J16").Value = 0.03654368891224535
$ws.Range("M16").Value = 89.83562999999999
$ws.Range("N16").Value = 269.50689
$ws.Range("O16").Value = 0.1681521229683693
$ws.Range("P16").Value = 0.1782820942415013
$ws.Range("Q16").Value = 5209.21147635848
$ws.Range("R16").Value = 46882.90328722631
$ws.Range("S16").Value = 0.006107076861403203
$ws.Range("T16").Value = 0.006515085390585031
$ws.Range("G17").Value = 29.481085
$ws.Range("H17").Value = 58.96217
$ws.Range("I17").Value = 0.01846507700595112
$ws.Range("J17").Value = 0.01238628926567028
$ws.Range("M17").Value = 121.928739
$ws.Range("N17").Value = 365.786217
$ws.Range("O17").Value = 0.2282232151508951
$ws.Range("P17").Value = 0.2419720431319445
$ws.Range("Q17").Value = 3594.591518401815
$ws.Range("R17").Value = 21567.54911041089
$ws.Range("S17").Value = 0.004214159242307029
$ws.Range("T17").Value = 0.00299713572043751
$ws.Range("G18").Value = 29.481085
$ws.Range("H18").Value = 58.96217
$ws.Range("I18").Value = 0.01846507700595112
$ws.Range("J18").Value = 0.01238628926567028
$ws.Range("O18").Value = 0.2768624053389947
$ws.Range("P18").Value = 0.2935413991166814
$ws.Range("Q18").Value = 4360.67493544805
$ws.Range("R18").Value = 26164.0496126883
$ws.Range("S18").Value = 0.005112285634637391
$ws.Range("T18").Value = 0.003635888680908786
$ws.Range("G19").Value = 29.481085
$ws.Range("H19").Value = 58.96217
$ws.Range("I19").Value = 0.01846507700595112
$ws.Range("J19").Value = 0.01238628926567028
$ws.Range("M19").Value = 83.50496933333334
$ws.Range("N19").Value = 250.514908
$ws.Range("O19").Value = 0.1563025480180701
$ws.Range("P19").Value = 0.1657186665504434
$ws.Range("Q19").Value = 2461.817098838394
$ws.Range("R19").Value = 14770.90259303036
$ws.Range("S19").Value = 0.002886138585380037
$ws.Range("T19").Value = 0.002052639340614949
$ws.Range("G20").Value = 29.481085
$ws.Range("H20").Value = 58.96217
$ws.Range("I20").Value = 0.01846507700595112
$ws.Range("J20").Value = 0.01238628926567028
$ws.Range("M20").Value = 91.06846250000001
$ws.Range("N20").Value = 182.136925
$ws.Range("O20").Value = 0.1704597085236707
$ws.Range("P20").Value = 0.1204857969594293
$ws.Range("Q20").Value = 2684.797083781813
$ws.Range("R20").Value = 10739.18833512725
$ws.Range("S20").Value = 0.003147551644301563
$ws.Range("T20").Value = 0.001492371933544308
$ws.Range("G21").Value = 29.481085
$ws.Range("H21").Value = 58.96217
$ws.Range("I21").Value = 0.01846507700595112
$ws.Range("J21").Value = 0.01238628926567028
$ws.Range("M21").Value = 89.83562999999999
$ws.Range("N21").Value = 269.50689
$ws.Range("O21").Value = 0.1681521229683693
$ws.Range("P21").Value = 0.1782820942415013
$ws.Range("Q21").Value = 2648.45184405855
$ws.Range("R21").Value = 15890.7110643513
$ws.Range("S21").Value = 0.003104941899325102
$ws.Range("T21").Value = 0.002208253590164724
$ws.Range("G22").Value = 21.628479
$ws.Range("H22").Value = 64.885437
$ws.Range("I22").Value = 0.01354670393768061
$ws.Range("J22").Value = 0.01363060063446486
$ws.Range("M22").Value = 121.928739
$ws.Range("N22").Value = 365.786217
$ws.Range("O22").Value = 0.2282232151508951
$ws.Range("P22").Value = 0.2419720431319445
$ws.Range("Q22").Value = 2637.133170957981
$ws.Range("R22").Value = 23734.19853862183
$ws.Range("S22").Value = 0.003091672327354759
$ws.Range("T22").Value = 0.003298224284637042
$ws.Range("G23").Value = 21.628479
$ws.Range("H23").Value = 64.885437
$ws.Range("I23").Value = 0.01354670393768061
$ws.Range("J23").Value = 0.01363060063446486
$ws.Range("O23").Value = 0.2768624053389947
$ws.Range("P23").Value = 0.2935413991166814
$ws.Range("Q23").Value = 3199.16198020407
$ws.Range("R23").Value = 28792.45782183663
$ws.Range("S23").Value = 0.003750573036601485
$ws.Range("T23").Value = 0.004001145581041541
$ws.Range("G24").Value = 21.628479
$ws.Range("H24").Value = 64.885437
$ws.Range("I24").Value = 0.01354670393768061
$ws.Range("J24").Value = 0.01363060063446486
$ws.Range("M24").Value = 83.50496933333334
$ws.Range("N24").Value = 250.514908
$ws.Range("O24").Value = 0.1563025480180701
$ws.Range("P24").Value = 0.1657186665504434
$ws.Range("Q24").Value = 1806.085475621644
$ws.Range("R24").Value = 16254.76928059479
$ws.Range("S24").Value = 0.002117384342705902
$ws.Range("T24").Value = 0.002258844961425144
$ws.Range("G25").Value = 21.628479
$ws.Range("H25").Value = 64.885437
$ws.Range("I25").Value = 0.01354670393768061
$ws.Range("J25").Value = 0.01363060063446486
$ws.Range("M25").Value = 91.06846250000001
$ws.Range("N25").Value = 182.136925
$ws.Range("O25").Value = 0.1704597085236707
$ws.Range("P25").Value = 0.1204857969594293
$ws.Range("Q25").Value = 1969.672328743537
$ws.Range("R25").Value = 11818.03397246123
$ws.Range("S25").Value = 0.002309167204673499
$ws.Range("T25").Value = 0.001642293780479202
$ws.Range("G26").Value = 21.628479
$ws.Range("H26").Value = 64.885437
$ws.Range("I26").Value = 0.01354670393768061
$ws.Range("J26").Value = 0.01363060063446486
$ws.Range("M26").Value = 89.83562999999999
$ws.Range("N26").Value = 269.50689
$ws.Range("O26").Value = 0.1681521229683693
$ws.Range("P26").Value = 0.1782820942415013
$ws.Range("Q26").Value = 1943.00803690677
$ws.Range("R26").Value = 17487.07233216093
$ws.Range("S26").Value = 0.002277907026344962
$ws.Range("T26").Value = 0.002430092026881932